# Re-theme the deck ("Integral" / Red Violet -> "Office Theme" / Office)
# and restyle the three data tables to the built-in "No Style, No Grid"
# table style.
#
# The presentation ships with two theme parts: the slide master / main
# design theme (Red Violet / "Integral") and the notes-master theme
# ("Office Theme"). The edit swaps their content so the design that
# drives the slides becomes the plain Office palette. The PowerPoint
# object model surfaces a single design/theme (the one driving the
# slide master), so we repaint its twelve theme colors to the Office
# palette values here.

$p = $ppt.ActivePresentation

# --- 1. Swap the design theme's colour scheme to the plain "Office" palette ---
$master = $p.Slides.Item(1).Master
$colors = $master.Theme.ThemeColorScheme

$colors.Colors(1).RGB  = 0        # dk1      000000
$colors.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388  # dk2      44546A
$colors.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501  # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407    # accent4  FFC000
$colors.Colors(9).RGB  = 12874308 # accent5  4472C4
$colors.Colors(10).RGB = 4697456  # accent6  70AD47
$colors.Colors(11).RGB = 12673797 # hlink    0563C1
$colors.Colors(12).RGB = 7491477  # folHlink 954F72

# --- 2. Re-apply table style "{9E7E2EE6-490D-4ACF-A61C-C5F7CABE0B72}" ---
#         (built-in "No Style, No Grid") to every table in the deck.
For ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    For ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        If ($shape.HasTable) {
            $shape.Table.ApplyStyle("{9E7E2EE6-490D-4ACF-A61C-C5F7CABE0B72}")
        }
    }
}
